# RF001 - Autenticar Usuario: bump test-suite wording from v1.3 to v1.4
#
# Changes:
#  - "inicia a tela de login" -> "abre a tela de login" (step 1 of every TC)
#  - For each test case (TC1..TC6) the "retry" step (row 3 of the 3-step
#    block) now performs the *other* user action (fill form vs. pick a
#    suggested user) relative to the one that failed in row 2, and the
#    failure message shown in row 2 is swapped accordingly between the
#    test cases so the overall transition-pair coverage stays balanced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abre = "Usuario do Sistema abre a tela de login atraves da opcao de Login no canto superior direito"
$preenche = "Usuario do Sistema preenche os campos e clica no botao entrar"
$seleciona = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$casFora = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"
$tjsegFora = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$senhaIncorreta = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"

# TC1 (rows 10-13)
$ws.Range("B10").Value = $abre
$ws.Range("D11").Value = $tjsegFora
$ws.Range("B12").Value = $seleciona

# TC2 (rows 20-23)
$ws.Range("B20").Value = $abre
$ws.Range("D21").Value = $casFora
$ws.Range("B22").Value = $preenche

# TC3 (rows 30-33)
$ws.Range("B30").Value = $abre
$ws.Range("B31").Value = $seleciona

# TC4 (rows 40-43)
$ws.Range("B40").Value = $abre
$ws.Range("B41").Value = $preenche
$ws.Range("D41").Value = $senhaIncorreta

# TC5 (rows 50-53)
$ws.Range("B50").Value = $abre
$ws.Range("D51").Value = $tjsegFora
$ws.Range("B52").Value = $preenche

# TC6 (rows 60-63)
$ws.Range("B60").Value = $abre
$ws.Range("D61").Value = $casFora
$ws.Range("B62").Value = $seleciona
